$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Category labels keyed by index (matches the recurring pattern already present
# in rows 2-14 of the sheet: Gehalt, Miete, Rewe, Kino, Rewe, Werkstatt,
# Restaurant, Rewe, Spotify, Rewe, Netflix, Gehalt, Miete, ...)
$categories = @("Gehalt","Miete","Rewe","Kino","Rewe","Werkstatt","Restaurant","Rewe","Spotify","Rewe","Netflix","Gehalt","Miete")
$amounts    = @(3000,-800,-100,-30,-250,-250,-100,-80,-10,-100,-10,3200,-800)

# Continue the repeating 13-row pattern for rows 15 through 39.
for ($row = 15; $row -le 39; $row++) {
    $idx = ($row - 2) % 13
    $ws.Cells.Item($row, 2).Value = $categories[$idx]
    $ws.Cells.Item($row, 3).Value = $amounts[$idx]
}

# Update the view state: scroll position and active selection.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("A10").Select()
$ws.Range("G37").Select()
